$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of SkillCode -> full Skill Description name
$skillNames = @{
    "SLEN" = "Systems and software life cycle engineering"
}

# Find the last used row based on column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Insert a new column before column B (shifts SFIA Level, Keycode, Description right)
$ws.Columns.Item(2).Insert(-4121)  # xlShiftToRight = -4121

# Set the new header
$ws.Cells.Item(1, 2).Value = "Skill Description"

# Populate the new column for each data row using the SkillCode in column A
for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($skillNames.ContainsKey($code)) {
        $ws.Cells.Item($r, 2).Value = $skillNames[$code]
    } else {
        $ws.Cells.Item($r, 2).Value = $code
    }
}
